$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New rows typed into Sheet2 (picked up into the shared-string table)
$ws2.Range("C2").Value = "jygkhjkhyk"
$ws2.Range("C7").Value = "desegfdsgf"

# Leave the cursor on Sheet1 at A3 before switching away from it
$ws1.Select()
$ws1.Range("A3").Select()

# Sheet2 becomes the active / selected tab, with C7 as the active cell
$ws2.Activate()
$ws2.Range("C7").Select()
